# Generate Report for Handback
# Update the recorded handoff/handback timestamps that are refreshed each
# time the handback status report is regenerated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for c5076bd4-...md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-07 05:38:16"

# --- zh-cn sheet: Correspond Handoff/Handback Datetime for c5076bd4-...zh-cn.xlf ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-07 05:37:59"
$wsZhCn.Range("K3").Value = "2016-09-07 05:38:50"

# --- de-de sheet: Correspond Handback Datetime for c5076bd4-...de-de.xlf ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-09-07 05:39:13"
